$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.852.47"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3
$ws.Range("D3").Value = "2.526.15"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("E9").Value = "  -1.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.64%  "

# Row 13
$ws.Range("E13").Value = "  -0.60%  "

# Row 14
$ws.Range("D14").Value = "2.914.93"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15
$ws.Range("D15").Value = "2.515.44"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.11%  "

# Row 17
$ws.Range("E17").Value = "  -1.68%  "

# Row 18
$ws.Range("D18").Value = "42.933.00"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20
$ws.Range("E20").Value = "  +3.27%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  -0.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("E24").Value = "  +1.26%  "

# Row 25
$ws.Range("E25").Value = "  -0.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  +3.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "

# Row 34
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0789"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.25%  "

# Row 37
$ws.Range("E37").Value = "  -0.82%  "

# Row 38
$ws.Range("E38").Value = "  -3.40%  "

# Row 39
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0303"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.99%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.07%  "

# Row 45
$ws.Range("E45").Value = "  -3.14%  "

# Row 46
$ws.Range("D46").Value = "2.018.27"
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.73%  "

# Row 49
$ws.Range("D49").Value = "2.769.73"
$ws.Range("E49").Value = "  +0.10%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
